$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EJ44")
$ws.Range("C3").Value = 226.1
$ws.Range("C4").Value = 233.2
$ws.Range("C5").Value = 238.9
$ws.Range("C6").Value = 243.8
$ws.Range("C7").Value = 248.5
$ws.Range("C8").Value = 253.6
$ws.Activate()
$ws.Range("C9").Select()
